$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Header date: "Curriculum Vitae, May 20, 2023" -> "..., August 1, 2023"
#    Use the unique preceding anchor "Curriculum Vitae, " so we do not
#    collide with other "May 20xx" occurrences later in the document.
# ------------------------------------------------------------------
$r1 = $d.Range(0, $d.Content.End)
$ok1 = $r1.Find.Execute("Curriculum Vitae, May 20, 2023", $true, $false, $false, $false, $false, $true, 1, $false, "Curriculum Vitae, August 1, 2023", 2)
Write-Host "date replace: " $ok1

# ------------------------------------------------------------------
# 2) Cisco Systems Inc. bullet: reword two phrases inside the paragraph.
# ------------------------------------------------------------------
$r2 = $d.Range(0, $d.Content.End)
$ok2 = $r2.Find.Execute("for the Cisco VDB team.", $true, $false, $false, $false, $false, $true, 1, $false, "for the Cisco Secure Firewall team.", 2)
Write-Host "cisco team replace: " $ok2

$r3 = $d.Range(0, $d.Content.End)
$ok3 = $r3.Find.Execute("Using my Big Data Analytics skillsets to identify patterns in network traffic to provide improved security for clients of", $true, $false, $false, $false, $false, $true, 1, $false, "Performing Big Data Analytics to identify patterns in network traffic to improve network security for clients of", 2)
Write-Host "cisco analytics replace: " $ok3

Write-Host "done"
